# The workbook's single sheet contains a list of municipalities (A),
# case counts (B) and death counts (C). Row 326 held a spurious
# "nafo informado" entry that needs to be removed entirely; deleting
# the row shifts every following row up by one, which also reduces
# the used range from A1:C573 to A1:C572.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 326 ("nafo informado"), shifting rows below it up.
$ws.Rows.Item(326).Delete()
